# "Corrige duplicado de issues"
#
# Row 3 is a continuation row for the same Host as row 2 (127.0.0.1), but
# its Host cell (B3) was left blank. Downstream tooling treated the blank
# Host as a separate/duplicate issue grouping. Fix: fill B3 with the same
# Host value as B2/B4, and leave the selection on the two Host cells that
# now match (B2:B3), mirroring what happened in Excel when making this fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing duplicated Host value for row 3.
$ws.Range("B3").Value = "127.0.0.1"

# Reflect the selection left behind after making the correction.
$ws.Range("B2:B3").Select()
